# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values on row 3 of the
# zh-cn and de-de worksheets to reflect the regenerated report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-02-15 03:48:28"
$wsZhCn.Range("G3").Value = "2016-02-15 03:49:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-02-15 03:48:42"
$wsDeDe.Range("G3").Value = "2016-02-15 03:49:39"
